$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ltp")

$values = @{
  2  = @(2506.45, 2472.45)
  3  = @(389.65, 391.6)
  4  = @(1503.45, 1499.9)
  5  = @(7293.35, 7235.6)
  6  = @(235.8, 234.9)
  7  = @(191.15, 190.15)
  8  = @(44500.5, 44486.65)
  9  = @(494, 482.4)
  10 = @(3355.35, 3348.35)
  11 = @(144.45, 143.9)
  12 = @(1157.55, 1146.45)
  13 = @(1400.15, 1409.45)
  14 = @(676.3, 655.85)
  15 = @(423.3, 420)
  16 = @(1560.85, 1549.25)
  17 = @(293.15, 292.7)
  18 = @(19337.1, 19321.3)
  19 = @(574.85, 572.85)
  20 = @(604.2, 602.4)
  21 = @(607.15, 602.4)
  22 = @(246.9, 244.95)
  23 = @(119.55, 117.6)
}

foreach ($row in $values.Keys) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 2).Value = $pair[0]
  $ws.Cells.Item($row, 3).Value = $pair[1]
}
